$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$halfPi = 1.5707963267948966

for ($row = 1; $row -le 17; $row++) {
    $bCell = $ws.Cells.Item($row, 2)  # Column B
    $cCell = $ws.Cells.Item($row, 3)  # Column C
    $dCell = $ws.Cells.Item($row, 4)  # Column D
    $eCell = $ws.Cells.Item($row, 5)  # Column E
    $fCell = $ws.Cells.Item($row, 6)  # Column F

    $bCell.Value = -1 * $bCell.Value2
    $cCell.Value = $halfPi - $cCell.Value2
    $dCell.Value = -1 * $dCell.Value2
    $eCell.Value = $halfPi - $eCell.Value2
    $fCell.Value = -1 * $fCell.Value2
}
